$p = $ppt.ActivePresentation

# --- Slide 2 ("Sprint 3 Stories"): "Total 68 pts" -> "Total 64 pts" ---
# TextBox 5 currently holds a single run "Total 68 pts". The score value
# changes from 68 to 64; editing just that portion of the text (as PowerPoint
# does when you select "68 " and retype it) splits the run into
# "Total " / "64 " / "pts".
$s2 = $p.Slides.Item(2)
$totalBox68 = $s2.Shapes.Item("TextBox 5")
$tr68 = $totalBox68.TextFrame.TextRange
$scorePart = $tr68.Characters(7, 3)
$scorePart.Text = "64 "

# --- Slide 11 ("Sprint 3 Stories"): "Total " + "36 pts" -> "Total 36 pts" ---
# TextBox 7 currently holds the phrase split across two runs ("Total " and
# "36 pts"). Retyping the whole phrase in one go collapses it back into a
# single run while keeping the original formatting.
$s11 = $p.Slides.Item(11)
$totalBox36 = $s11.Shapes.Item("TextBox 7")
$tr36 = $totalBox36.TextFrame.TextRange
$whole = $tr36.Characters(1, $tr36.Text.Length)
$whole.Text = "Total 36 pts"
